$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" header columns (A:J) to "_FV2404" and the
# "_new" header columns (L:U) to "_FV2410". Column K holds "diff"
# and is left untouched.
$oldLetters = @("A","B","C","D","E","F","G","H","I","J")
$newLetters = @("L","M","N","O","P","Q","R","S","T","U")
$baseNames  = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($oldLetters[$i] + "1").Value = $baseNames[$i] + "_FV2404"
    $ws.Range($newLetters[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a table with an autofilter on the header row.
$rng = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
